$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price values in column D stay as text (matches source inlineStr cells)
$priceCells = @("D2", "D3", "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D34", "D38", "D40", "D41", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = '57.218.45'
$ws.Range("E2").Value = '  -0.73%  '
$ws.Range("D3").Value = '2.408.11'
$ws.Range("E3").Value = '  -1.55%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '512.12'
$ws.Range("E5").Value = '  -2.18%  '
$ws.Range("D6").Value = '133.21'
$ws.Range("E6").Value = '  +2.15%  '
$ws.Range("D7").Value = '0.993'
$ws.Range("E7").Value = '  -0.53%  '
$ws.Range("D8").Value = '0.559'
$ws.Range("E8").Value = '  -0.51%  '
$ws.Range("D9").Value = '2.445.76'
$ws.Range("E9").Value = '  +0.09%  '
$ws.Range("D10").Value = '0.0975'
$ws.Range("E10").Value = '  +0.02%  '
$ws.Range("E11").Value = '  -0.78%  '
$ws.Range("D12").Value = '0.323'
$ws.Range("E12").Value = '  +0.45%  '
$ws.Range("D13").Value = '4.71'
$ws.Range("E13").Value = '  -4.43%  '
$ws.Range("D14").Value = '2.833.63'
$ws.Range("E14").Value = '  -1.82%  '
$ws.Range("D15").Value = '56.996.53'
$ws.Range("E15").Value = '  -1.05%  '
$ws.Range("D16").Value = '21.95'
$ws.Range("E16").Value = '  +1.13%  '
$ws.Range("D17").Value = '0.0000134'
$ws.Range("E17").Value = '  +0.98%  '
$ws.Range("D18").Value = '2.378.84'
$ws.Range("E18").Value = '  -2.79%  '
$ws.Range("D19").Value = '10.30'
$ws.Range("E19").Value = '  -0.38%  '
$ws.Range("B20").Value = 'Polkadot'
$ws.Range("C20").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D20").Value = '4.13'
$ws.Range("E20").Value = '  -0.24%  '
$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").Value = '313.95'
$ws.Range("E21").Value = '  -0.26%  '
$ws.Range("D22").Value = '6.45'
$ws.Range("E22").Value = '  +6.54%  '
$ws.Range("D23").Value = '0.997'
$ws.Range("E23").Value = '  -0.20%  '
$ws.Range("D24").Value = '5.72'
$ws.Range("E24").Value = '  -1.94%  '
$ws.Range("D25").Value = '65.12'
$ws.Range("E25").Value = '  +0.70%  '
$ws.Range("D26").Value = '0.990'
$ws.Range("E26").Value = '  -1.41%  '
$ws.Range("D27").Value = '2.490.83'
$ws.Range("E27").Value = '  -3.08%  '
$ws.Range("D28").Value = '0.383'
$ws.Range("E28").Value = '  -5.89%  '
$ws.Range("D29").Value = '0.151'
$ws.Range("E29").Value = '  -3.58%  '
$ws.Range("D30").Value = '7.50'
$ws.Range("E30").Value = '  +3.38%  '
$ws.Range("D31").Value = '172.84'
$ws.Range("E31").Value = '  +0.23%  '
$ws.Range("D32").Value = '0.0₃0740'
$ws.Range("E32").Value = '  +0.77%  '
$ws.Range("E33").Value = '  -0.05%  '
$ws.Range("D34").Value = '6.23'
$ws.Range("E34").Value = '  +2.12%  '
$ws.Range("E35").Value = '  -0.09%  '
$ws.Range("E37").Value = '  -0.59%  '
$ws.Range("D38").Value = '18.07'
$ws.Range("E38").Value = '  +1.70%  '
$ws.Range("E39").Value = '  +4.24%  '
$ws.Range("D40").Value = '3.86'
$ws.Range("E40").Value = '  +1.61%  '
$ws.Range("D41").Value = '0.821'
$ws.Range("E41").Value = '  +3.64%  '
$ws.Range("E42").Value = '  -0.19%  '
$ws.Range("D43").Value = '36.06'
$ws.Range("E43").Value = '  -0.31%  '
$ws.Range("B44").Value = 'Filecoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D44").Value = '3.43'
$ws.Range("E44").Value = '  +0.52%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").Value = '131.14'
$ws.Range("E45").Value = '  +5.61%  '
$ws.Range("D46").Value = '4.99'
$ws.Range("E46").Value = '  +3.52%  '
$ws.Range("D47").Value = '260.38'
$ws.Range("E47").Value = '  -2.12%  '
$ws.Range("D48").Value = '0.571'
$ws.Range("E48").Value = '  -2.00%  '
$ws.Range("D49").Value = '0.0912'
$ws.Range("E49").Value = '  -1.66%  '
$ws.Range("D50").Value = '0.0497'
$ws.Range("E50").Value = '  +0.55%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '17.37'
$ws.Range("E51").Value = '  +2.09%  '
